$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H82").Value = 24643.334
$ws_ALC.Range("I82").Value = 1941.5
$ws_ALC.Range("J82").Value = 70047
$ws_ALC.Range("K82").Value = 5824.5
$ws_ALC.Range("L82").Value = 210141
$ws_ALC.Range("M82").Value = -5418.5
$ws_ALC.Range("N82").Value = -210953

$ws_ALC.Range("H85").Value = 24643.334
$ws_ALC.Range("I85").Value = 1941.5
$ws_ALC.Range("J85").Value = 70047
$ws_ALC.Range("K85").Value = 5824.5
$ws_ALC.Range("L85").Value = 210141
$ws_ALC.Range("M85").Value = -4420.5
$ws_ALC.Range("N85").Value = -212949

$ws_ALC.Range("H98").Value = 52747852
$ws_ALC.Range("I98").Value = 83342430
$ws_ALC.Range("J98").Value = 300000
$ws_ALC.Range("K98").Value = 83342430
$ws_ALC.Range("L98").Value = 300000
$ws_ALC.Range("M98").Value = -83340932

$ws_ALC.Range("H122").Value = 52747852
$ws_ALC.Range("I122").Value = 83342430
$ws_ALC.Range("J122").Value = 300000
$ws_ALC.Range("K122").Value = 250027290
$ws_ALC.Range("L122").Value = 900000
$ws_ALC.Range("M122").Value = -250024840

$ws_ALC.Range("H138").Value = 2986.5396
$ws_ALC.Range("I138").Value = 951.4
$ws_ALC.Range("J138").Value = 3294.894
$ws_ALC.Range("K138").Value = 2854.2
$ws_ALC.Range("L138").Value = 9884.681999999999
$ws_ALC.Range("M138").Value = 2285.8
$ws_ALC.Range("N138").Value = -20164.682

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H5").Value = 313.4
$ws_ARM.Range("I5").Value = 313.4
$ws_ARM.Range("J5").Value = 0
$ws_ARM.Range("K5").Value = 313.4
$ws_ARM.Range("L5").Value = 0
$ws_ARM.Range("M5").Value = -201.4
$ws_ARM.Range("N5").ClearContents()

$ws_ARM.Range("H61").Value = 31321102
$ws_ARM.Range("I61").Value = 83340536
$ws_ARM.Range("J61").Value = 109438.1
$ws_ARM.Range("K61").Value = 83340536
$ws_ARM.Range("L61").Value = 109438.1
$ws_ARM.Range("M61").Value = -83340324
$ws_ARM.Range("N61").Value = -109862.1

$ws_ARM.Range("H136").Value = 31321102
$ws_ARM.Range("I136").Value = 83340536
$ws_ARM.Range("J136").Value = 109438.1
$ws_ARM.Range("K136").Value = 250021608
$ws_ARM.Range("L136").Value = 328314.3
$ws_ARM.Range("M136").Value = -250019058
$ws_ARM.Range("N136").Value = -333414.3

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H4").Value = 313.4
$ws_BSM.Range("I4").Value = 313.4
$ws_BSM.Range("J4").Value = 0
$ws_BSM.Range("K4").Value = 313.4
$ws_BSM.Range("L4").Value = 0
$ws_BSM.Range("M4").Value = -198.4
$ws_BSM.Range("N4").ClearContents()

$ws_BSM.Range("H107").Value = 0
$ws_BSM.Range("I107").Value = 0
$ws_BSM.Range("J107").Value = 0
$ws_BSM.Range("K107").Value = 0
$ws_BSM.Range("L107").Value = 0
$ws_BSM.Range("M107").ClearContents()

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H7").Value = 2957.4167
$ws_CRP.Range("I7").Value = 95
$ws_CRP.Range("J7").Value = 6964.8
$ws_CRP.Range("K7").Value = 95
$ws_CRP.Range("L7").Value = 6964.8
$ws_CRP.Range("M7").Value = 18
$ws_CRP.Range("N7").Value = -7190.8

$ws_CRP.Range("H19").Value = 900
$ws_CRP.Range("I19").Value = 900
$ws_CRP.Range("J19").Value = 0
$ws_CRP.Range("K19").Value = 900
$ws_CRP.Range("L19").Value = 0
$ws_CRP.Range("M19").Value = -730
$ws_CRP.Range("N19").ClearContents()

$ws_CRP.Range("H24").Value = 900
$ws_CRP.Range("I24").Value = 900
$ws_CRP.Range("J24").Value = 0
$ws_CRP.Range("K24").Value = 900
$ws_CRP.Range("L24").Value = 0
$ws_CRP.Range("M24").Value = -730
$ws_CRP.Range("N24").ClearContents()

$ws_CRP.Range("H58").Value = 2125.5
$ws_CRP.Range("I58").Value = 2197.8125
$ws_CRP.Range("J58").Value = 1547
$ws_CRP.Range("K58").Value = 2197.8125
$ws_CRP.Range("L58").Value = 1547
$ws_CRP.Range("M58").Value = -1994.8125
$ws_CRP.Range("N58").Value = -1953

$ws_CRP.Range("H96").Value = 6099
$ws_CRP.Range("I96").Value = 1999
$ws_CRP.Range("J96").Value = 6919
$ws_CRP.Range("K96").Value = 1999
$ws_CRP.Range("L96").Value = 6919
$ws_CRP.Range("M96").Value = 747
$ws_CRP.Range("N96").Value = -12411

$ws_CRP.Range("H107").Value = 1728
$ws_CRP.Range("I107").Value = 1179.6
$ws_CRP.Range("J107").Value = 2642
$ws_CRP.Range("K107").Value = 1179.6
$ws_CRP.Range("L107").Value = 2642
$ws_CRP.Range("M107").Value = 740.4000000000001

$ws_CRP.Range("H132").Value = 2384.8572
$ws_CRP.Range("I132").Value = 2384.8572
$ws_CRP.Range("J132").Value = 0
$ws_CRP.Range("K132").Value = 7154.571599999999
$ws_CRP.Range("L132").Value = 0
$ws_CRP.Range("M132").Value = -4624.571599999999

$ws_CRP.Range("H134").Value = 198734.47
$ws_CRP.Range("I134").Value = 257968.17
$ws_CRP.Range("J134").Value = 6224.9165
$ws_CRP.Range("K134").Value = 773904.51
$ws_CRP.Range("L134").Value = 18674.7495
$ws_CRP.Range("M134").Value = -771369.51
$ws_CRP.Range("N134").Value = -23744.7495

$ws_CRP.Range("H136").Value = 2125.5
$ws_CRP.Range("I136").Value = 2197.8125
$ws_CRP.Range("J136").Value = 1547
$ws_CRP.Range("K136").Value = 6593.4375
$ws_CRP.Range("L136").Value = 4641
$ws_CRP.Range("M136").Value = -4043.4375
$ws_CRP.Range("N136").Value = -9741

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 712.4666999999999
$ws_CUL.Range("I5").Value = 720.5
$ws_CUL.Range("J5").Value = 600
$ws_CUL.Range("K5").Value = 2161.5
$ws_CUL.Range("L5").Value = 1800
$ws_CUL.Range("M5").Value = -2049.5
$ws_CUL.Range("N5").Value = -2024

$ws_CUL.Range("H34").Value = 829.8
$ws_CUL.Range("I34").Value = 537.25
$ws_CUL.Range("J34").Value = 2000
$ws_CUL.Range("K34").Value = 1611.75
$ws_CUL.Range("L34").Value = 6000
$ws_CUL.Range("M34").Value = -1527.75
$ws_CUL.Range("N34").Value = -6168

$ws_CUL.Range("H131").Value = 5903.3335
$ws_CUL.Range("I131").Value = 7068.875
$ws_CUL.Range("J131").Value = 4571.2856
$ws_CUL.Range("K131").Value = 21206.625
$ws_CUL.Range("L131").Value = 13713.8568
$ws_CUL.Range("M131").Value = -16166.625
$ws_CUL.Range("N131").Value = -23793.8568

$ws_CUL.Range("H135").Value = 712.4666999999999
$ws_CUL.Range("I135").Value = 720.5
$ws_CUL.Range("J135").Value = 600
$ws_CUL.Range("K135").Value = 6484.5
$ws_CUL.Range("L135").Value = 5400
$ws_CUL.Range("M135").Value = -3949.5
$ws_CUL.Range("N135").Value = -10470

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H2").Value = 236.93333
$ws_GSM.Range("I2").Value = 32.6
$ws_GSM.Range("J2").Value = 645.6
$ws_GSM.Range("K2").Value = 32.6
$ws_GSM.Range("L2").Value = 645.6
$ws_GSM.Range("M2").Value = 80.40000000000001

$ws_GSM.Range("H57").Value = 9000.6
$ws_GSM.Range("I57").Value = 9000.6
$ws_GSM.Range("J57").Value = 0
$ws_GSM.Range("K57").Value = 9000.6
$ws_GSM.Range("L57").Value = 0
$ws_GSM.Range("M57").Value = -8180.6

$ws_GSM.Range("H132").Value = 100003320
$ws_GSM.Range("I132").Value = 125003400
$ws_GSM.Range("J132").Value = 3000
$ws_GSM.Range("K132").Value = 375010200
$ws_GSM.Range("L132").Value = 9000
$ws_GSM.Range("M132").Value = -375007670

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H68").Value = 3671.3333
$ws_LTW.Range("I68").Value = 3498.5
$ws_LTW.Range("J68").Value = 3887.375
$ws_LTW.Range("K68").Value = 3498.5
$ws_LTW.Range("L68").Value = 3887.375
$ws_LTW.Range("M68").Value = -2749.5
$ws_LTW.Range("N68").Value = -5385.375

$ws_LTW.Range("H71").Value = 3671.3333
$ws_LTW.Range("I71").Value = 3498.5
$ws_LTW.Range("J71").Value = 3887.375
$ws_LTW.Range("K71").Value = 17492.5
$ws_LTW.Range("L71").Value = 19436.875
$ws_LTW.Range("M71").Value = -13748.5
$ws_LTW.Range("N71").Value = -26924.875

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H64").Value = 64996.668
$ws_WVR.Range("I64").Value = 0
$ws_WVR.Range("J64").Value = 64996.668
$ws_WVR.Range("K64").Value = 0
$ws_WVR.Range("L64").Value = 64996.668
$ws_WVR.Range("N64").Value = -65492.668
$ws_WVR.Range("M64").ClearContents()

$ws_WVR.Range("H67").Value = 64996.668
$ws_WVR.Range("I67").Value = 0
$ws_WVR.Range("J67").Value = 64996.668
$ws_WVR.Range("K67").Value = 0
$ws_WVR.Range("L67").Value = 64996.668
$ws_WVR.Range("N67").Value = -66712.66800000001
$ws_WVR.Range("M67").ClearContents()

$ws_WVR.Range("H107").Value = 26316800
$ws_WVR.Range("I107").Value = 29412812
$ws_WVR.Range("J107").Value = 690
$ws_WVR.Range("K107").Value = 88238436
$ws_WVR.Range("L107").Value = 2070
$ws_WVR.Range("M107").Value = -88236516
